$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for column AY (30/09/2021) across the balance-sheet rows.
$updates = @{
    2  = 12955074.56
    3  = 2768220.928
    4  = 1043358.016
    5  = 89943
    6  = 638864
    7  = 239000.992
    9  = 157676.992
    10 = 0
    11 = 599377.9840000001
    12 = 2408293.888
    13 = 0
    14 = 289
    19 = 53541
    20 = 0
    22 = 0
    23 = 5968074.24
    24 = 1810486.016
    26 = 12955074.56
    27 = 11214551.04
    28 = 373416.992
    29 = 1773283.968
    30 = 61876
    31 = 1848112
    34 = 6827213.824
    35 = 330648
    37 = 20020137.984
    38 = 9950574.592
    40 = 8302731.776
    41 = 10211
    43 = 1756620.032
    46 = 0
    47 = -18279614.464
    48 = 4041435.904
    49 = 160744
    51 = 0
    52 = -21397583.872
    53 = -1085273.984
    55 = 1065
    59 = 2922333.952
    60 = -3873302.528
    61 = -950969.088
    62 = -228828.96
    63 = -603430.0159999999
    65 = -15678
    66 = 49516
    68 = -1092075.008
    69 = -446390.08
    70 = -645685.12
    74 = -2841464.832
    75 = 82
    76 = 32060
    79 = 0
    80 = -2809322.752
}

foreach ($row in $updates.Keys) {
    $ws.Range("AY$row").Value = $updates[$row]
}
